# ---------------------------------------------------------------------------
# Reproduces the authored edit:
#   1. Slide 16's table switches to table style {CD0D1B71-85A5-49B7-89E3-4EA7ED12087D}
#      (was {5DB562D3-1EFB-4CA6-862C-5F3C1CA5E717}).
#   2. The theme used by the slide master / presentation (ppt/theme/theme2.xml,
#      the "Integral" theme) is recoloured to the stock "Office Theme" palette
#      (the palette that theme1.xml - the notes-master theme - already carries).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style -----------------------------------------------------
$slide16 = $p.Slides.Item(16)
for ($i = 1; $i -le $slide16.Shapes.Count; $i++) {
    $shp = $slide16.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{CD0D1B71-85A5-49B7-89E3-4EA7ED12087D}")
    }
}

# --- 2. Theme colours -----------------------------------------------------
# Slide 1's ThemeColorScheme is backed by the one theme part that the slide
# master (and therefore every slide) draws from.
$tcs = $p.Slides.Item(1).ThemeColorScheme

function Set-ThemeRGB {
    param($scheme, [int]$index, [byte]$r, [byte]$g, [byte]$b)
    $scheme.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

Set-ThemeRGB $tcs 1  0x00 0x00 0x00   # dk1
Set-ThemeRGB $tcs 2  0xFF 0xFF 0xFF   # lt1
Set-ThemeRGB $tcs 3  0x44 0x54 0x6A   # dk2
Set-ThemeRGB $tcs 4  0xE7 0xE6 0xE6   # lt2
Set-ThemeRGB $tcs 5  0x5B 0x9B 0xD5   # accent1
Set-ThemeRGB $tcs 6  0xED 0x7D 0x31   # accent2
Set-ThemeRGB $tcs 7  0xA5 0xA5 0xA5   # accent3
Set-ThemeRGB $tcs 8  0xFF 0xC0 0x00   # accent4
Set-ThemeRGB $tcs 9  0x44 0x72 0xC4   # accent5
Set-ThemeRGB $tcs 10 0x70 0xAD 0x47   # accent6
Set-ThemeRGB $tcs 11 0x05 0x63 0xC1   # hlink
Set-ThemeRGB $tcs 12 0x95 0x4F 0x72   # folHlink
